$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Handback report generation: the localization pipeline has produced a
# handback for the e498de5f-... source file in both zh-cn and de-de, so the
# "Ready for handoff" status everywhere becomes "Handed back: in sync with
# en-US", and the per-language detail sheets get their Target/Handback file
# columns and Handback DateTime filled in.
# ---------------------------------------------------------------------------

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$mdDisplay = "e498de5f-2453-4d31-81c0-8ce5f05474c0.md"
$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/fada21acb90c9c3ccbaee13e2456b884eab16dc8/e2e/e498de5f-2453-4d31-81c0-8ce5f05474c0.md"

# ---------------------------------------------------------------------------
# Overview sheet: Status columns (B = zh-cn, C = de-de) for both rows.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value2 = $statusNew
$wsOverview.Range("C2").Value2 = $statusNew
$wsOverview.Range("B3").Value2 = $statusNew
$wsOverview.Range("C3").Value2 = $statusNew

# ---------------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$zhXlfDisplay = "e498de5f-2453-4d31-81c0-8ce5f05474c0.88694f32d0f0406db45b4acd55aaf8271b2f35be.zh-cn.xlf"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/70ac6b3fd81b421f5a2e1d09ce3cc5eb65688994/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e498de5f-2453-4d31-81c0-8ce5f05474c0.88694f32d0f0406db45b4acd55aaf8271b2f35be.zh-cn.xlf"
$zhHandbackTime = "2016-03-23 23:11:24"

$wsZh.Range("C2").Value2 = $statusNew
$wsZh.Range("C3").Value2 = $statusNew

$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdDisplay) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $zhXlfUrl, [Type]::Missing, [Type]::Missing, $zhXlfDisplay) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdDisplay) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $zhXlfUrl, [Type]::Missing, [Type]::Missing, $zhXlfDisplay) | Out-Null

$wsZh.Range("H2").Value2 = $zhHandbackTime
$wsZh.Range("H3").Value2 = $zhHandbackTime

# ---------------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$deXlfDisplay = "e498de5f-2453-4d31-81c0-8ce5f05474c0.88694f32d0f0406db45b4acd55aaf8271b2f35be.de-de.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/40baab0a5a2d3612d7a512a1f26143cf885c96fc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e498de5f-2453-4d31-81c0-8ce5f05474c0.88694f32d0f0406db45b4acd55aaf8271b2f35be.de-de.xlf"
$deHandbackTime = "2016-03-23 23:11:31"

$wsDe.Range("C2").Value2 = $statusNew
$wsDe.Range("C3").Value2 = $statusNew

$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdDisplay) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $deXlfUrl, [Type]::Missing, [Type]::Missing, $deXlfDisplay) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdDisplay) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $deXlfUrl, [Type]::Missing, [Type]::Missing, $deXlfDisplay) | Out-Null

$wsDe.Range("H2").Value2 = $deHandbackTime
$wsDe.Range("H3").Value2 = $deHandbackTime
